$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "output" rows (MCU_REDUNDANCY_1/2 + new FAN_ENABLE) into the
#     H:L output table at rows 8-10 ---
$ws.Range("H8").Value = "MCU_REDUNDANCY_1"
$ws.Range("I8").Value = "digital"
$ws.Range("J8").Value = "PTD0"
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 2

$ws.Range("H9").Value = "MCU_REDUNDANCY_2"
$ws.Range("I9").Value = "digital"
$ws.Range("J9").Value = "PTD1"
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1

$ws.Range("H10").Value = "FAN_ENABLE"
$ws.Range("I10").Value = "digital"
$ws.Range("J10").Value = "PTD2"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 46

# --- Remove the now-duplicated MCU_REDUNDANCY_1/2 rows from the bottom
#     standalone table (B:F), shifting the remaining rows (CAN GENERAL
#     TX/RX, pwr, signal pwr, pwr gnd, sig gnd) up by two ---
$ws.Range("B21:F22").Delete(-4162) | Out-Null

# --- Update the selected cell to match the edited workbook's view state ---
$ws.Range("L11").Select() | Out-Null
